$d = $word.ActiveDocument

$d.Content.Find.Execute("89×12=1068", $true, $false, $false, $false, $false, $true, 1, $false, "21×26=546", 2) | Out-Null
$d.Content.Find.Execute("43×31=1333", $true, $false, $false, $false, $false, $true, 1, $false, "37×18=666", 2) | Out-Null
$d.Content.Find.Execute("13×85=1105", $true, $false, $false, $false, $false, $true, 1, $false, "45×15=675", 2) | Out-Null
$d.Content.Find.Execute("57×42=2394", $true, $false, $false, $false, $false, $true, 1, $false, "15×41=615", 2) | Out-Null
$d.Content.Find.Execute("73×37=2701", $true, $false, $false, $false, $false, $true, 1, $false, "16×34=544", 2) | Out-Null
$d.Content.Find.Execute("71×70=4970", $true, $false, $false, $false, $false, $true, 1, $false, "94×18=1692", 2) | Out-Null
$d.Content.Find.Execute("52×49=2548", $true, $false, $false, $false, $false, $true, 1, $false, "78×74=5772", 2) | Out-Null
$d.Content.Find.Execute("74×76=5624", $true, $false, $false, $false, $false, $true, 1, $false, "87×45=3915", 2) | Out-Null
$d.Content.Find.Execute("76×87=6612", $true, $false, $false, $false, $false, $true, 1, $false, "63×20=1260", 2) | Out-Null
$d.Content.Find.Execute("42×17=714", $true, $false, $false, $false, $false, $true, 1, $false, "42×96=4032", 2) | Out-Null
$d.Content.Find.Execute("13×34=442", $true, $false, $false, $false, $false, $true, 1, $false, "98×67=6566", 2) | Out-Null
$d.Content.Find.Execute("60×79=4740", $true, $false, $false, $false, $false, $true, 1, $false, "50×42=2100", 2) | Out-Null
$d.Content.Find.Execute("29×39=1131", $true, $false, $false, $false, $false, $true, 1, $false, "80×19=1520", 2) | Out-Null
$d.Content.Find.Execute("18×72=1296", $true, $false, $false, $false, $false, $true, 1, $false, "46×46=2116", 2) | Out-Null
$d.Content.Find.Execute("46×12=552", $true, $false, $false, $false, $false, $true, 1, $false, "11×47=517", 2) | Out-Null
$d.Content.Find.Execute("22×43=946", $true, $false, $false, $false, $false, $true, 1, $false, "55×56=3080", 2) | Out-Null
$d.Content.Find.Execute("61×98=5978", $true, $false, $false, $false, $false, $true, 1, $false, "26×23=598", 2) | Out-Null
$d.Content.Find.Execute("12×71=852", $true, $false, $false, $false, $false, $true, 1, $false, "41×81=3321", 2) | Out-Null
$d.Content.Find.Execute("29×57=1653", $true, $false, $false, $false, $false, $true, 1, $false, "67×23=1541", 2) | Out-Null
$d.Content.Find.Execute("99×12=1188", $true, $false, $false, $false, $false, $true, 1, $false, "31×38=1178", 2) | Out-Null
$d.Content.Find.Execute("35×25=875", $true, $false, $false, $false, $false, $true, 1, $false, "85×98=8330", 2) | Out-Null
$d.Content.Find.Execute("69×42=2898", $true, $false, $false, $false, $false, $true, 1, $false, "54×34=1836", 2) | Out-Null
$d.Content.Find.Execute("82×98=8036", $true, $false, $false, $false, $false, $true, 1, $false, "26×23=598", 2) | Out-Null
$d.Content.Find.Execute("11×93=1023", $true, $false, $false, $false, $false, $true, 1, $false, "38×26=988", 2) | Out-Null
$d.Content.Find.Execute("48×44=2112", $true, $false, $false, $false, $false, $true, 1, $false, "39×29=1131", 2) | Out-Null
